$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F21").Value = 3906
$ws1.Range("F31").Value = 581
$ws1.Range("F34").Value = 968
$ws1.Range("F35").Value = 2476

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F21").Value = 3906
$ws4.Range("F31").Value = 581
$ws4.Range("F35").Value = 968
$ws4.Range("F36").Value = 2476
